# Weekly update: insert a new data row for "Terminal La Palmera de La Serena - Jengibre"
# at row 33 (pushing the existing rows 33-39 down to 34-40), matching the
# commit "Fruta / hortaliza, semanal".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 33; existing rows 33:39 shift down to 34:40 and the
# style of row 33 (incl. the date-format style on column D) is carried down
# with them, so the new row 33 needs its own values written in.
$ws.Rows.Item(33).Insert()

$ws.Cells.Item(33, 1).Value = 8
$ws.Cells.Item(33, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(33, 3).Value = "Coquimbo"
$ws.Cells.Item(33, 4).Value = 44694
$ws.Cells.Item(33, 5).Value = 4
$ws.Cells.Item(33, 6).Value = 100114007
$ws.Cells.Item(33, 7).Value = "Jengibre"
$ws.Cells.Item(33, 8).Value = "Sin especificar"
$ws.Cells.Item(33, 9).Value = "Primera"
$ws.Cells.Item(33, 10).Value = 400
$ws.Cells.Item(33, 11).Value = 13000
$ws.Cells.Item(33, 12).Value = 14000
$ws.Cells.Item(33, 13).Value = 13500
$ws.Cells.Item(33, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(33, 15).Value = "Perú"
$ws.Cells.Item(33, 16).Value = 1038
$ws.Cells.Item(33, 17).Value = 13
$ws.Cells.Item(33, 18).Value = "Hortaliza"
